$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddr, $value)
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "68.191.95"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "2.516.40"
$ws.Range("E3").Value = "  +1.17%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "589.91"
$ws.Range("E5").Value = "  +1.02%  "
Set-TextValue "D6" "177.65"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("E12").Value = "  +0.62%  "
Set-TextValue "D13" "25.86"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "2.851.06"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").Value = "68.010.04"
$ws.Range("E15").Value = "  +1.61%  "
Set-TextValue "D16" "0.0000172"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "2.503.22"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +2.12%  "
Set-TextValue "D20" "352.66"
$ws.Range("E20").Value = "  +1.14%  "
Set-TextValue "D21" "4.06"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.18%  "
Set-TextValue "D23" "70.72"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "2.642.77"
$ws.Range("E27").Value = "  +1.06%  "
Set-TextValue "D28" "0.998"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  +1.50%  "
Set-TextValue "D30" "510.98"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E35").Value = "  +4.32%  "
Set-TextValue "D36" "165.05"
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +3.67%  "
Set-TextValue "D43" "0.330"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  +5.17%  "
Set-TextValue "D45" "147.33"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +3.65%  "
Set-TextValue "D49" "0.0744"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +2.43%  "
Set-TextValue "D51" "0.589"
$ws.Range("E51").Value = "  +1.02%  "
